$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($Cell, $Text) {
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.Style = "Normal"
}

Set-TextCell $ws.Cells.Item(2, 4) "29.187.27"
Set-TextCell $ws.Cells.Item(2, 5) "  -0.65%  "

Set-TextCell $ws.Cells.Item(3, 4) "1.857.90"
Set-TextCell $ws.Cells.Item(3, 5) "  -1.25%  "

Set-TextCell $ws.Cells.Item(4, 5) "  -0.03%  "

Set-TextCell $ws.Cells.Item(5, 4) "242.13"
Set-TextCell $ws.Cells.Item(5, 5) "  -0.10%  "

Set-TextCell $ws.Cells.Item(6, 4) "0.7024"
Set-TextCell $ws.Cells.Item(6, 5) "  -1.53%  "

Set-TextCell $ws.Cells.Item(7, 5) "  +0.00%  "

Set-TextCell $ws.Cells.Item(8, 4) "0.3113"
Set-TextCell $ws.Cells.Item(8, 5) "  -0.45%  "

Set-TextCell $ws.Cells.Item(9, 4) "0.07777"
Set-TextCell $ws.Cells.Item(9, 5) "  -3.45%  "

Set-TextCell $ws.Cells.Item(10, 5) "  -4.52%  "

Set-TextCell $ws.Cells.Item(11, 4) "0.07982"
Set-TextCell $ws.Cells.Item(11, 5) "  -4.17%  "

Set-TextCell $ws.Cells.Item(12, 4) "1.856.33"
Set-TextCell $ws.Cells.Item(12, 5) "  -1.83%  "

Set-TextCell $ws.Cells.Item(13, 4) "93.58"
Set-TextCell $ws.Cells.Item(13, 5) "  -0.18%  "

Set-TextCell $ws.Cells.Item(14, 4) "5.165"
Set-TextCell $ws.Cells.Item(14, 5) "  -1.54%  "

Set-TextCell $ws.Cells.Item(15, 4) "0.6950"
Set-TextCell $ws.Cells.Item(15, 5) "  -3.32%  "

Set-TextCell $ws.Cells.Item(16, 4) "6.365"
Set-TextCell $ws.Cells.Item(16, 5) "  +0.57%  "

Set-TextCell $ws.Cells.Item(17, 4) "29.178.94"
Set-TextCell $ws.Cells.Item(17, 5) "  -0.72%  "

Set-TextCell $ws.Cells.Item(18, 4) "0.000008277"
Set-TextCell $ws.Cells.Item(18, 5) "  -3.52%  "

Set-TextCell $ws.Cells.Item(19, 4) "250.99"
Set-TextCell $ws.Cells.Item(19, 5) "  +3.69%  "

Set-TextCell $ws.Cells.Item(20, 4) "2.109.25"
Set-TextCell $ws.Cells.Item(20, 5) "  -1.46%  "

Set-TextCell $ws.Cells.Item(21, 4) "13.08"
Set-TextCell $ws.Cells.Item(21, 5) "  -1.21%  "

Set-TextCell $ws.Cells.Item(22, 5) "  +0.00%  "

Set-TextCell $ws.Cells.Item(23, 4) "7.500"
Set-TextCell $ws.Cells.Item(23, 5) "  -4.59%  "

Set-TextCell $ws.Cells.Item(24, 5) "  +0.01%  "

Set-TextCell $ws.Cells.Item(25, 4) "0.1552"
Set-TextCell $ws.Cells.Item(25, 5) "  -2.44%  "

Set-TextCell $ws.Cells.Item(26, 4) "8.973"
Set-TextCell $ws.Cells.Item(26, 5) "  -1.04%  "

Set-TextCell $ws.Cells.Item(27, 4) "159.27"
Set-TextCell $ws.Cells.Item(27, 5) "  -2.69%  "

Set-TextCell $ws.Cells.Item(28, 4) "18.83"
Set-TextCell $ws.Cells.Item(28, 5) "  +1.21%  "

Set-TextCell $ws.Cells.Item(29, 4) "1.496"
Set-TextCell $ws.Cells.Item(29, 5) "  -0.67%  "

Set-TextCell $ws.Cells.Item(30, 4) "4.296"
Set-TextCell $ws.Cells.Item(30, 5) "  -2.65%  "

Set-TextCell $ws.Cells.Item(31, 4) "4.258"
Set-TextCell $ws.Cells.Item(31, 5) "  -1.84%  "

Set-TextCell $ws.Cells.Item(32, 4) "1.210"
Set-TextCell $ws.Cells.Item(32, 5) "  +0.87%  "

Set-TextCell $ws.Cells.Item(33, 4) "0.05247"
Set-TextCell $ws.Cells.Item(33, 5) "  -2.27%  "

Set-TextCell $ws.Cells.Item(34, 4) "1.873"
Set-TextCell $ws.Cells.Item(34, 5) "  -3.87%  "

Set-TextCell $ws.Cells.Item(35, 4) "0.7416"
Set-TextCell $ws.Cells.Item(35, 5) "  -1.02%  "

Set-TextCell $ws.Cells.Item(36, 5) "  -2.52%  "

Set-TextCell $ws.Cells.Item(37, 4) "2.712"
Set-TextCell $ws.Cells.Item(37, 5) "  +0.63%  "

Set-TextCell $ws.Cells.Item(38, 4) "0.01865"
Set-TextCell $ws.Cells.Item(38, 5) "  -1.18%  "

Set-TextCell $ws.Cells.Item(39, 4) "1.241.01"
Set-TextCell $ws.Cells.Item(39, 5) "  -3.57%  "

Set-TextCell $ws.Cells.Item(40, 5) "  -0.43%  "

Set-TextCell $ws.Cells.Item(41, 4) "6.160"
Set-TextCell $ws.Cells.Item(41, 5) "  -6.69%  "

Set-TextCell $ws.Cells.Item(42, 2) "Quant"
Set-TextCell $ws.Cells.Item(42, 3) "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextCell $ws.Cells.Item(42, 4) "110.61"
Set-TextCell $ws.Cells.Item(42, 5) "  -1.01%  "

Set-TextCell $ws.Cells.Item(43, 2) "TrustWalletToken"
Set-TextCell $ws.Cells.Item(43, 3) "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell $ws.Cells.Item(43, 4) "0.8932"
Set-TextCell $ws.Cells.Item(43, 5) "  -2.58%  "

Set-TextCell $ws.Cells.Item(44, 4) "70.74"
Set-TextCell $ws.Cells.Item(44, 5) "  -5.48%  "

Set-TextCell $ws.Cells.Item(45, 5) "  -0.04%  "

Set-TextCell $ws.Cells.Item(46, 5) "  -0.04%  "

Set-TextCell $ws.Cells.Item(47, 4) "2.007.45"
Set-TextCell $ws.Cells.Item(47, 5) "  -1.29%  "

Set-TextCell $ws.Cells.Item(48, 4) "0.5178"
Set-TextCell $ws.Cells.Item(48, 5) "  -0.82%  "

Set-TextCell $ws.Cells.Item(49, 4) "1.778"
Set-TextCell $ws.Cells.Item(49, 5) "  -1.77%  "

Set-TextCell $ws.Cells.Item(50, 4) "9.407"
Set-TextCell $ws.Cells.Item(50, 5) "  -1.09%  "

Set-TextCell $ws.Cells.Item(51, 4) "0.4291"
Set-TextCell $ws.Cells.Item(51, 5) "  -2.20%  "
